$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: add label "Time" in A12, update F12 and G12 values
$ws.Range("A12").Value = "Time"
$ws.Range("F12").Value = 6.66
$ws.Range("G12").Value = 3.15

# Row 13: F13 keeps its existing formula (recalculates automatically from
# the new F12 value); G13 gets its own new formula (no longer shared with F13)
$ws.Range("G13").Formula = "=G12/`$F12"

# Update the active selection to F13
$ws.Range("F13").Select()
